$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the header text in A1 (shared string previously
# "A adm pública tem agido nas causas LGBTQIA+?") with "prefeitura".
$ws.Range("A1").Value = "prefeitura"

# Update the active selection to D3:E3 (active cell D3).
[void]$ws.Range("D3:E3").Select()
